# Update "想去人数" (want-to-go count) figures for two events that changed
# between scrape runs:
#   - 南宁·熊喵M动漫嘉年华·万圣派对   132 -> 140
#   - 南宁·万圣漫控嘉年华10           688 -> 693
#
# These events appear on the "展览" sheet (rows 3 & 4) and are duplicated
# on the "全部类型" sheet (rows 4 & 5) since it aggregates all events.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 140
$wsExhibit.Range("F4").Value = 693

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 140
$wsAll.Range("F5").Value = 693
